# Journal de travail - Daniel : "Started writing report. Updated jdt"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Row 50 was "Travail en groupe via chat vocal" (no hours logged yet).
# Update the activity description to reflect the work actually done and
# log the 0.5 hour spent on it.
$ws.Range("B50").Value = "Rédaction d'une partie du rapport concernant la logique métier."
$ws.Range("C50").Value = 0.5

# Move the active selection down to A51, where the next entry will go.
$ws.Range("A51").Select()

$wb.Save()
